# Adds the "Abstract Title" paragraph style, tightens the spacing above
# the existing "Abstract" style, and adds the "Footnote Block Text"
# paragraph style (based on "Footnote Text"), mirroring the upstream
# styles.xml diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New "Abstract Title" style -- centered, bold, small, blue heading
#    that precedes the "Abstract" style paragraph.
# ---------------------------------------------------------------------

$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# ---------------------------------------------------------------------
# 2. Existing "Abstract" style: pull its leading space in from 15pt
#    (300 twips) to 5pt (100 twips); trailing space is unchanged.
# ---------------------------------------------------------------------

$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------
# 3. New "Footnote Block Text" style, based on "Footnote Text", used for
#    block-quoted material inside footnotes.
# ---------------------------------------------------------------------

$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "styles updated"
